$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue "D2" '67.289.29'
Set-TextValue "E2" '  +1.48%  '

# Row 3
Set-TextValue "D3" '3.872.86'
Set-TextValue "E3" '  +1.02%  '

# Row 4
Set-TextValue "E4" '  -0.02%  '

# Row 5
Set-TextValue "D5" '468.33'
Set-TextValue "E5" '  +9.40%  '

# Row 6
Set-TextValue "D6" '144.95'
Set-TextValue "E6" '  +10.37%  '

# Row 7
Set-TextValue "D7" '0.633'

# Row 8
Set-TextValue "D8" '0.998'
Set-TextValue "E8" '  -0.12%  '

# Row 9
Set-TextValue "D9" '0.746'
Set-TextValue "E9" '  +1.98%  '

# Row 10
Set-TextValue "E10" '  -1.10%  '

# Row 11
Set-TextValue "D11" '0.0000312'
Set-TextValue "E11" '  -6.30%  '

# Row 12
Set-TextValue "D12" '43.35'
Set-TextValue "E12" '  +3.81%  '

# Row 13
Set-TextValue "D13" '10.45'
Set-TextValue "E13" '  -0.35%  '

# Row 14
Set-TextValue "D14" '4.492.18'
Set-TextValue "E14" '  +0.94%  '

# Row 15
Set-TextValue "D15" '14.81'
Set-TextValue "E15" '  -5.43%  '

# Row 16
Set-TextValue "D16" '3.874.75'
Set-TextValue "E16" '  +0.25%  '

# Row 17
Set-TextValue "E17" '  -0.38%  '

# Row 18
Set-TextValue "D18" '20.05'
Set-TextValue "E18" '  -0.10%  '

# Row 19
Set-TextValue "E19" '  +6.00%  '

# Row 20
Set-TextValue "D20" '67.539.55'
Set-TextValue "E20" '  +1.41%  '

# Row 21
Set-TextValue "D21" '436.50'
Set-TextValue "E21" '  +4.89%  '

# Row 22
Set-TextValue "D22" '14.89'
Set-TextValue "E22" '  -1.05%  '

# Row 23
Set-TextValue "E23" '  +6.31%  '

# Row 24
Set-TextValue "D24" '89.12'
Set-TextValue "E24" '  +4.68%  '

# Row 25
Set-TextValue "D25" '3.59'
Set-TextValue "E25" '  +8.76%  '

# Row 26
Set-TextValue "D26" '38.06'
Set-TextValue "E26" '  +1.80%  '

# Row 27
Set-TextValue "D27" '10.09'
Set-TextValue "E27" '  +8.12%  '

# Row 28
Set-TextValue "D28" '9.92'
Set-TextValue "E28" '  -1.73%  '

# Row 29
Set-TextValue "D29" '5.55'
Set-TextValue "E29" '  +2.93%  '

# Row 30
Set-TextValue "D30" '729.08'
Set-TextValue "E30" '  +1.08%  '

# Row 31
Set-TextValue "D31" '13.84'
Set-TextValue "E31" '  -0.93%  '

# Row 32
Set-TextValue "E32" '  +6.49%  '

# Row 33
Set-TextValue "E33" '  +0.50%  '

# Row 34
Set-TextValue "D34" '44.20'
Set-TextValue "E34" '  +13.20%  '

# Row 35
Set-TextValue "D35" '0.162'
Set-TextValue "E35" '  +8.38%  '

# Row 36
Set-TextValue "D36" '58.18'
Set-TextValue "E36" '  +4.48%  '

# Row 37
Set-TextValue "B37" 'Dai'
Set-TextValue "C37" 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue "D37" '1.00'
Set-TextValue "E37" '  +0.08%  '

# Row 38
Set-TextValue "B38" 'NEARProtocol'
Set-TextValue "C38" 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D38" '5.50'
Set-TextValue "E38" '  -5.08%  '

# Row 39
Set-TextValue "D39" '0.0484'
Set-TextValue "E39" '  +2.61%  '

# Row 40
Set-TextValue "D40" '2.92'
Set-TextValue "E40" '  +0.97%  '

# Row 41
Set-TextValue "D41" '0.346'

# Row 42
Set-TextValue "B42" 'Stellar'
Set-TextValue "C42" 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue "D42" '0.142'
Set-TextValue "E42" '  +3.65%  '

# Row 43
Set-TextValue "B43" 'PEPE'
Set-TextValue "C43" 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue "D43" '0.0₃0679'
Set-TextValue "E43" '  -6.89%  '

# Row 44
Set-TextValue "E44" '  -0.06%  '

# Row 45
Set-TextValue "E45" '  +5.69%  '

# Row 46
Set-TextValue "D46" '3.46'
Set-TextValue "E46" '  +1.66%  '

# Row 47
Set-TextValue "D47" '3.28'
Set-TextValue "E47" '  -0.02%  '

# Row 48
Set-TextValue "E48" '  +5.08%  '

# Row 49
Set-TextValue "D49" '2.15'
Set-TextValue "E49" '  +4.57%  '

# Row 50
Set-TextValue "B50" 'Monero'
Set-TextValue "C50" 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue "D50" '144.58'
Set-TextValue "E50" '  +1.68%  '

# Row 51
Set-TextValue "B51" 'Stacks'
Set-TextValue "C51" 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue "D51" '2.90'
Set-TextValue "E51" '  +1.44%  '
